$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-16 20:42:11"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-16 20:41:58"
$wsZhCn.Range("K4").Value = "2016-08-16 20:42:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-16 20:42:11"
$wsDeDe.Range("K4").Value = "2016-08-16 20:42:37"
